$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Altan)
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 5.333333333333333
$ws.Range("J2").Value = $true
$ws.Range("M2").Value = 1
$ws.Range("O2").Value = 1

# Row 3 (Goromi)
$ws.Range("D3").Value = 10
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 21
$ws.Range("H3").Value = 7
$ws.Range("I3").Value = $true
$ws.Range("L3").Value = 1
$ws.Range("O3").Value = 2

# Row 4 (Cyber)
$ws.Range("D4").Value = 9
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 6.333333333333333
$ws.Range("P4").Value = 1

# Row 5 (Nozomi)
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 12
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = $false
$ws.Range("L5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0

# Row 6 (Inizio)
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = 13
$ws.Range("H6").Value = 4.333333333333333
$ws.Range("P6").Value = 0
